# ============================================================================
# feat: add 2022-Q1 data
#
# - Inserts a new "2022-Q1" worksheet (fund-holding detail, same layout as
#   the existing quarterly sheets) right before the "总计" (totals) sheet.
# - Inserts a new leading row into the "总计" sheet summarising 2022-Q1
#   (21 funds, 7.25 亿元 held), shifting the prior rows down and bumping
#   their running index in column A.
# ============================================================================

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# 1) Insert the new "2022-Q1" worksheet right after "2021-Q4" (so it ends
#    up directly before "总计", matching the target sheet order).
# ----------------------------------------------------------------------
$afterSheet = $wb.Worksheets.Item("2021-Q4")
$q1 = $wb.Worksheets.Add($null, $afterSheet)
$q1.Name = "2022-Q1"

# Re-use the "2021-Q4" sheet (identical report layout) purely as a
# formatting template for the header row and the A-column row-index style.
$template = $wb.Worksheets.Item("2021-Q4")

# Header row B1:H1 -> same bold / bordered / centered style used by every
# other quarterly sheet.
$template.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)

$q1.Cells.Item(1,2).Value = "基金代码"
$q1.Cells.Item(1,3).Value = "基金名称"
$q1.Cells.Item(1,4).Value = "基金规模"
$q1.Cells.Item(1,5).Value = "股票总仓位"
$q1.Cells.Item(1,6).Value = "仓位占比"
$q1.Cells.Item(1,7).Value = "持有市值(亿元)"
$q1.Cells.Item(1,8).Value = "仓位排名"

# Column A (row index) style, rows 2-22, copied from the template sheet.
$template.Range("A2").Copy()
$q1.Range("A2:A22").PasteSpecial(-4122)

# ----------------------------------------------------------------------
# Fund-holding detail rows (A..H), row 2 .. row 22.
#   A: running index (number)
#   B: fund code          (text - keep leading zeros, e.g. "010676")
#   C: fund name           (text)
#   D: fund scale (亿元)   (text, matches source formatting)
#   E: total equity position (text)
#   F: position weight     (text)
#   G: holding value (亿元) (text)
#   H: position rank        (number)
# ----------------------------------------------------------------------
$q1Data = @(
    ,@(0, "206009", "鹏华新兴产业混合", "44.95", "90.17", "4.99", "2.2430", 6)
    ,@(1, "519068", "汇添富成长焦点混合", "54.15", "85.21", "3.22", "1.7436", 6)
    ,@(2, "310328", "申万菱信新动力混合", "34.56", "73.99", "2.69", "0.9297", 7)
    ,@(3, "360006", "光大保德信新增长混合", "21.71", "88.07", "4.25", "0.9227", 4)
    ,@(4, "310308", "申万菱信盛利精选混合", "14.72", "66.62", "2.60", "0.3827", 10)
    ,@(5, "010676", "光大保德信新机遇混合", "4.20", "85.57", "6.44", "0.2705", 1)
    ,@(6, "008878", "国联安新蓝筹红利一年定期开放混合", "3.96", "97.84", "4.67", "0.1849", 3)
    ,@(7, "009794", "太平智选一年定期开放股票", "5.20", "89.96", "2.60", "0.1352", 10)
    ,@(8, "360005", "光大保德信红利混合", "4.34", "85.30", "2.66", "0.1154", 9)
    ,@(9, "006568", "国联安行业领先混合", "1.58", "91.33", "6.52", "0.1030", 3)
    ,@(10, "005708", "国联安远见成长混合", "1.52", "91.59", "6.42", "0.0976", 3)
    ,@(11, "001412", "德邦鑫星价值灵活配置混合A", "1.83", "68.16", "4.25", "0.0778", 4)
    ,@(12, "007903", "长城量化小盘股票", "1.36", "90.57", "1.17", "0.0159", 4)
    ,@(13, "011800", "申万菱信价值精选混合型证券投资基金", "0.57", "81.46", "2.70", "0.0154", 7)
    ,@(14, "002112", "德邦鑫星价值灵活配置混合C", "0.15", "68.16", "4.25", "0.0064", 4)
    ,@(15, "004726", "先锋聚优灵活配置混合A", "0.06", "93.23", "5.47", "0.0033", 2)
    ,@(16, "004727", "先锋聚优灵活配置混合C", "0.04", "93.23", "5.47", "0.0022", 2)
    ,@(17, "003586", "先锋精一灵活配置混合A", "0.03", "92.66", "3.69", "0.0011", 10)
    ,@(18, "003587", "先锋精一灵活配置混合C", "0.03", "92.66", "3.69", "0.0011", 10)
    ,@(19, "004833", "先锋聚利灵活配置混合A", "0.02", "94.68", "4.78", "0.0010", 6)
    ,@(20, "004834", "先锋聚利灵活配置混合C", "0.02", "94.68", "4.78", "0.0010", 6)
)

# Force columns B:G to Text storage *before* writing, so numeric-looking
# strings (fund codes with leading zeros, "44.95", etc.) are not silently
# coerced to numbers. Clear the number-format residue afterwards so the
# cells end up with no explicit style, matching the other quarterly sheets.
$q1.Range("B2:G22").NumberFormat = "@"

$r = 2
foreach ($row in $q1Data) {
    $q1.Cells.Item($r,1).Value = $row[0]
    $q1.Cells.Item($r,2).Value = $row[1]
    $q1.Cells.Item($r,3).Value = $row[2]
    $q1.Cells.Item($r,4).Value = $row[3]
    $q1.Cells.Item($r,5).Value = $row[4]
    $q1.Cells.Item($r,6).Value = $row[5]
    $q1.Cells.Item($r,7).Value = $row[6]
    $q1.Cells.Item($r,8).Value = $row[7]
    $r = $r + 1
}

$q1.Range("B2:G22").ClearFormats()

# ----------------------------------------------------------------------
# 2) "总计" sheet: insert a new leading data row (row 2) summarising
#    2022-Q1, pushing the previous rows down by one.
# ----------------------------------------------------------------------
$totals = $wb.Worksheets.Item("总计")
$totals.Rows.Item(2).Insert()

# Give the new row's A-cell the same style as the rest of column A, and
# clear the inherited header-ish formatting the row Insert() leaves on
# B2:D2 so it matches the plain (unstyled) data rows below it.
$totals.Range("A3").Copy()
$totals.Range("A2").PasteSpecial(-4122)
$totals.Range("B2:D2").ClearFormats()

$totals.Cells.Item(2,1).Value = 0
$totals.Cells.Item(2,2).Value = "2022-Q1"
$totals.Cells.Item(2,3).Value = 21
$totals.Cells.Item(2,4).Value = 7.25

# Bump the running index in column A for the rows that shifted down.
for ($row = 3; $row -le 7; $row++) {
    $totals.Cells.Item($row,1).Value = $row - 2
}
